# Allow excluding transaction's category in reports.
#
# On the "kategorikas" sheet, a new "dipakaiDiLaporan" (used in report)
# boolean column is inserted before the existing "sistem" column: the
# existing "sistem" header/values move from column I to the new column J,
# and column I becomes the new "dipakaiDiLaporan" column (re-using the
# old values, which were already TRUE for every data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kategorikas")

# --- move the existing "sistem" column (I) out to the new column (J) ---

$sistemHeader = $ws.Range("I1").Value()
$ws.Range("J1").Value = $sistemHeader
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").Interior.Color = $ws.Range("I1").Interior.Color()

$ws.Range("J2").Value = $ws.Range("I2").Value()
$ws.Range("J3").Value = $ws.Range("I3").Value()
$ws.Range("J4").Value = $ws.Range("I4").Value()
$ws.Range("J5").Value = $ws.Range("I5").Value()

# --- turn column I into the new "dipakaiDiLaporan" column ---

$ws.Range("I1").Value = "dipakaiDiLaporan"
$ws.Range("I2").Value = $true
$ws.Range("I3").Value = $true
$ws.Range("I4").Value = $true
$ws.Range("I5").Value = $true

# --- column widths: widen I for the longer header, restore J to the ---
# --- width the "sistem" column used to have                        ---

$ws.Columns.Item(9).ColumnWidth = 16.67
$ws.Columns.Item(10).ColumnWidth = 5.83

# --- make "kategorikas" the active sheet / tab, selection on I4 ---

$ws.Activate()
$ws.Range("I4").Select()
